$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.706.61"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "2.379.83"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("E7").Value = "  -4.34%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.113"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "2.751.94"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("D17").Value = "2.384.61"
$ws.Range("E17").Value = "  +3.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.817"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "43.633.35"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.87%  "
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").Value = "0.0₃0921"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.34%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +6.01%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.62%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  +6.88%  "
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.76%  "
$ws.Range("D45").Value = "2.038.58"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0291"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.55%  "
$ws.Range("E51").Value = "  +0.67%  "
